$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.315.97"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.680.21"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'218.25"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'0.5257"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "'0.2694"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").Value = "'0.06463"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").Value = "'21.99"
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("D11").Value = "'0.07522"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "1.689.28"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "'4.528"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'0.5811"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "'0.000008499"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "'64.83"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "26.348.46"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").Value = "'4.922"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").Value = "'189.89"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").Value = "'6.212"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'145.32"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'7.828"
$ws.Range("E25").Value = "  +2.72%  "
$ws.Range("D26").Value = "'0.1253"
$ws.Range("E26").Value = "  +3.87%  "
$ws.Range("D27").Value = "'15.81"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").Value = "'0.06479"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("D29").Value = "'1.361"
$ws.Range("E29").Value = "  +5.01%  "
$ws.Range("D30").Value = "'1.326"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "'3.608"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("D32").Value = "'3.597"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").Value = "'1.669"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("D34").Value = "'1.031"
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").Value = "'0.6247"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").Value = "'2.407"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("D38").Value = "'6.434"
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("D39").Value = "'0.01625"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").Value = "1.108.74"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("D41").Value = "'0.8777"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").Value = "'100.60"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "1.831.64"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'56.97"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000106"
$ws.Range("E46").Value = "  -7.37%  "
$ws.Range("D47").Value = "'8.235"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").Value = "'0.05271"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  +2.33%  "
